# Updates cryptos list price/volume/name/link values to match the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.209.24'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.377.81'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.02%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.84'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +9.08%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.377.70'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.03%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.60'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.53%  '
$ws.Range('E11').Value = '  +5.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.392'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.963.41'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.122'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000173'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.382.79'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.15'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.427.36'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.94'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '380.90'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.575'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.51%  '
$ws.Range('E24').Value = '  +2.23%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.60%  '
$ws.Range('E27').Value = '  +13.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +15.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.82'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +10.21%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.16'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.84%  '
$ws.Range('E32').Value = '  +7.27%  '
$ws.Range('E33').Value = '  +2.16%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.415.56'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.49'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.56'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.23%  '
$ws.Range('E38').Value = '  +5.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.55'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '162.59'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  +8.01%  '
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.43'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.54%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.48'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.761'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.55%  '
$ws.Range('E46').Value = '  +8.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.68'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +9.70%  '
$ws.Range('E48').Value = '  +4.03%  '
$ws.Range('E49').Value = '  +5.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.16'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +12.76%  '
$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.887'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.84%  '
